$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$target  = $ws.Cells.Item(11, 2)   # B11
$scratch = $ws.Cells.Item(20, 20)  # unused helper cell

# Excel's Range.Value setter auto-parses digit-only strings as numbers, so a
# plain `$target.Value = "1"` would store a numeric 1 instead of the text
# "1" the target shared-string table expects. Build the literal text "1" on
# a throwaway cell that's forced to Text format first ...
$scratch.NumberFormat = "@"
$scratch.Value = "1"

# ... then copy only the VALUE over to B11, leaving B11's own number
# format/style (border, fill, etc.) completely untouched.
$scratch.Copy()
$target.PasteSpecial(-4163)   # xlPasteValues

# Tidy up the scratch cell and clipboard marching ants.
$scratch.Clear()
$excel.CutCopyMode = $false
